$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-08 Friday" "2024-03-09 Saturday"

Replace-Text "395×8=3160" "315×6=1890"
Replace-Text "699×9=6291" "133×3=399"
Replace-Text "208×9=1872" "414×3=1242"
Replace-Text "125×5=625" "773×2=1546"
Replace-Text "125×9=1125" "951×7=6657"
Replace-Text "689×9=6201" "698×9=6282"
Replace-Text "193×2=386" "628×9=5652"
Replace-Text "655×6=3930" "640×6=3840"
Replace-Text "789×2=1578" "336×6=2016"
Replace-Text "648×6=3888" "368×9=3312"
Replace-Text "826×8=6608" "327×8=2616"
Replace-Text "187×8=1496" "124×8=992"
Replace-Text "771×8=6168" "840×8=6720"
Replace-Text "853×4=3412" "774×5=3870"
Replace-Text "654×6=3924" "342×9=3078"
Replace-Text "829×6=4974" "969×2=1938"
Replace-Text "382×4=1528" "596×5=2980"
Replace-Text "989×3=2967" "454×3=1362"
Replace-Text "617×7=4319" "314×4=1256"
Replace-Text "242×2=484" "937×3=2811"
Replace-Text "675×9=6075" "789×8=6312"
Replace-Text "452×7=3164" "411×8=3288"
Replace-Text "135×7=945" "503×9=4527"
Replace-Text "969×8=7752" "644×3=1932"
Replace-Text "243×4=972" "225×4=900"
